$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace numeric data in rows 2-5 (34 columns: A..AH) with the new dataset
$row2 = @(45056.50694444445, 15.374, 10.112, 3.717, 32.732, 25.126, 11.943, 36.363, 18.615, 7.558, 11.183, 12.932, 13.607, 3.86, 12.031, 16.608, 10.57, 3.118, 1.764, 175.75, 33.355, 11.105, 21.604, 11.537, 2.926, 18.421, 9.808999999999999, 8.874000000000001, 10.611, 13.63, 3.312, 32.661, 5.894, 13.883)
$row3 = @(45056.51388888889, 13.452, 9.414999999999999, 1.677, 29.127, 22.997, 10.461, 40.459, 16.288, 6.969, 10.149, 11.629, 12.291, 3.382, 10.527, 14.77, 9.202, 1.447, 0.9350000000000001, 152.911, 29.467, 9.717000000000001, 19.397, 10.372, 2.023, 19.571, 8.583, 7.761, 9.167999999999999, 12.229, 1.246, 37.16, 5.277, 12.148)
$row4 = @(45056.52083333334, 2.883, 1.712, 0.855, 6.214, 4.375, 2.164, 13.835, 3.49, 1.421, 1.784, 2.472, 2.63, 0.74, 2.256, 3.141, 2.229, 0.895, 0.41, 27.083, 6.651, 2.082, 4.206, 2.272, 0.724, 6.081, 1.839, 1.795, 2.111, 2.593, 0.766, 13.198, 0.998, 2.608)
$row5 = @(45056.52777777778, 0.96, 0.38, 0.58, 2.06, 1.11, 0.67, 5.94, 1.16, 0.42, 0.35, 0.82, 0.89, 0.26, 0.75, 1.01, 0.91, 0.66, 0.26, 4.17, 2.32, 0.6899999999999999, 1.35, 0.79, 0.42, 2.39, 0.61, 0.68, 0.79, 0.85, 0.55, 5.72, 0.25, 0.87)

$rowsData = @{
  2 = $row2
  3 = $row3
  4 = $row4
  5 = $row5
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# Row 6 (old extra observation) is removed entirely; dataset now spans rows 1-5
$ws.Rows("6:6").Delete()

# Widen a subset of data columns from 7 to 8 character units (custom accuracy columns)
$wideCols = @("C", "G", "K", "L", "O", "V", "X", "AC")
foreach ($col in $wideCols) {
    $ws.Columns($col + ":" + $col).ColumnWidth = 7.166666666666667
}
